$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 241
$ws.Range("I2").Value = 662
$ws.Range("J2").Value = 2644
$ws.Range("K2").Value = 13
$ws.Range("L2").Value = 688
$ws.Range("M2").Value = 39
$ws.Range("N2").Value = 446
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 33
$ws.Range("S2").Value = 313
$ws.Range("T2").Value = 469
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 4078
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 4085
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 66
